$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "formula" row (row 6: formula / Datatype that holds formula / let f1 = FORMULA([a, b], SUM);)
# This shifts the "class" row (old row 7) up into row 6, and all string indices
# shift accordingly.
$ws.Rows.Item(6).Delete()

# Update the remaining "Example" column cells that now use `new` before the
# constructor-style calls.
$ws.Range("C5").Value = "let tab1 = new TABLE();"
$ws.Range("C2").Value = "let A1 = new CELL(5);"
$ws.Range("C6").Value = "struct example {`r`n    let id = 1, speaker = ""john"";`r`n    let country = ""au"";`r`n};`r`nconst exampleObj = new example(101, ""widget"", 19.99);"

# Adjust row heights for rows 4 and 5 to match the updated layout.
$ws.Rows.Item(4).RowHeight = 21
$ws.Rows.Item(5).RowHeight = 32

# Update the active selection to reflect the new last-used cell.
$ws.Range("C6").Select()
